# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# Price (D) and Volume(1h) (E) columns are plain text in this sheet (e.g.
# "61.603.08" uses '.' as a thousands separator, not a valid number), so for
# any new D-column value that *would* parse as a plain number we force the
# cell to Text format first - otherwise Excel's automatic type detection
# would silently convert e.g. "1.00" -> the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.324.44'
$ws.Range('E2').Value = '  -3.82%  '

$ws.Range('D3').Value = '3.001.65'
$ws.Range('E3').Value = '  -3.02%  '

$ws.Range('E4').Value = '  +0.27%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '537.38'
$ws.Range('E5').Value = '  -0.87%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.00'
$ws.Range('E6').Value = '  -0.69%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = '2.999.97'
$ws.Range('E8').Value = '  -2.75%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -0.48%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.149'
$ws.Range('E10').Value = '  -4.46%  '

$ws.Range('E11').Value = '  -0.19%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.450'
$ws.Range('E12').Value = '  -2.25%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000222'
$ws.Range('E13').Value = '  -2.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.34'
$ws.Range('E14').Value = '  -1.51%  '

$ws.Range('D15').Value = '3.492.28'
$ws.Range('E15').Value = '  -2.89%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.110'
$ws.Range('E16').Value = '  -1.07%  '

$ws.Range('D17').Value = '61.499.43'
$ws.Range('E17').Value = '  -3.32%  '

$ws.Range('D18').Value = '3.005.71'
$ws.Range('E18').Value = '  -2.81%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.66'
$ws.Range('E19').Value = '  -0.99%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '468.50'
$ws.Range('E20').Value = '  -4.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.30'
$ws.Range('E21').Value = '  -1.59%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.678'
$ws.Range('E22').Value = '  -3.72%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.98'
$ws.Range('E23').Value = '  -3.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.08'
$ws.Range('E24').Value = '  +0.17%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.08'
$ws.Range('E25').Value = '  -1.63%  '

$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.69'
$ws.Range('E27').Value = '  -1.95%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.86'
$ws.Range('E28').Value = '  -5.98%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.90'
$ws.Range('E30').Value = '  -1.15%  '

$ws.Range('B31').Value = 'Mantle'
$ws.Range('C31').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +3.04%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.69'
$ws.Range('E32').Value = '  -2.51%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.51'
$ws.Range('E33').Value = '  +1.64%  '

$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '55.52'
$ws.Range('E34').Value = '  -3.06%  '

$ws.Range('B35').Value = 'Stacks'
$ws.Range('C35').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.29'
$ws.Range('E35').Value = '  -5.25%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.91'
$ws.Range('E36').Value = '  -2.93%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '455.69'
$ws.Range('E37').Value = '  -7.56%  '

$ws.Range('D38').Value = '3.177.38'
$ws.Range('E38').Value = '  -0.87%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0793'
$ws.Range('E39').Value = '  -1.51%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0387'
$ws.Range('E40').Value = '  -3.66%  '

$ws.Range('E41').Value = '  +0.29%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.17'
$ws.Range('E42').Value = '  -0.19%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.50'
$ws.Range('E43').Value = '  -7.96%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '27.55'
$ws.Range('E44').Value = '  +11.05%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.246'
$ws.Range('E46').Value = '  -4.19%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.01'
$ws.Range('E47').Value = '  -2.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '120.06'
$ws.Range('E48').Value = '  -1.08%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.108'
$ws.Range('E49').Value = '  -0.79%  '

$ws.Range('D50').Value = '0.0₃0499'
$ws.Range('E50').Value = '  -8.82%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.02'
$ws.Range('E51').Value = '  -0.96%  '
